$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.969508
$ws.Range("H2").Value = 116.908524
$ws.Range("I2").Value = 0.688733638790647
$ws.Range("J2").Value = 0.688733638790647
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7341896666666666
$ws.Range("N2").Value = 2.202569
$ws.Range("O2").Value = 0.4912907638668469
$ws.Range("P2").Value = 0.4912907638668469
$ws.Range("Q2").Value = 28.611010088684
$ws.Range("R2").Value = 257.499090798156
$ws.Range("S2").Value = 0.33836847550225
$ws.Range("T2").Value = 0.33836847550225
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.969508
$ws.Range("H3").Value = 116.908524
$ws.Range("I3").Value = 0.688733638790647
$ws.Range("J3").Value = 0.688733638790647
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7602199999999999
$ws.Range("N3").Value = 2.28066
$ws.Range("O3").Value = 0.5087092361331531
$ws.Range("P3").Value = 0.5087092361331531
$ws.Range("Q3").Value = 29.62539937175999
$ws.Range("R3").Value = 266.62859434584
$ws.Range("S3").Value = 0.350365163288397
$ws.Range("T3").Value = 0.350365163288397
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.685730333333334
$ws.Range("H4").Value = 29.057191
$ws.Range("I4").Value = 0.1711822560557247
$ws.Range("J4").Value = 0.1711822560557247
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7341896666666666
$ws.Range("N4").Value = 2.202569
$ws.Range("O4").Value = 0.4912907638668469
$ws.Range("P4").Value = 0.4912907638668469
$ws.Range("Q4").Value = 7.111163124853222
$ws.Range("R4").Value = 64.00046812367901
$ws.Range("S4").Value = 0.0841002613380672
$ws.Range("T4").Value = 0.08410026133806718
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.685730333333334
$ws.Range("H5").Value = 29.057191
$ws.Range("I5").Value = 0.1711822560557247
$ws.Range("J5").Value = 0.1711822560557247
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7602199999999999
$ws.Range("N5").Value = 2.28066
$ws.Range("O5").Value = 0.5087092361331531
$ws.Range("P5").Value = 0.5087092361331531
$ws.Range("Q5").Value = 7.363285914006666
$ws.Range("R5").Value = 66.26957322606
$ws.Range("S5").Value = 0.08708199471765755
$ws.Range("T5").Value = 0.08708199471765755
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6067633333333333
$ws.Range("H6").Value = 1.82029
$ws.Range("I6").Value = 0.01072372580252768
$ws.Range("J6").Value = 0.01072372580252768
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.7341896666666666
$ws.Range("N6").Value = 2.202569
$ws.Range("O6").Value = 0.4912907638668469
$ws.Range("P6").Value = 0.4912907638668469
$ws.Range("Q6").Value = 0.4454793694455556
$ws.Range("R6").Value = 4.00931432501
$ws.Range("S6").Value = 0.005268467441022441
$ws.Range("T6").Value = 0.005268467441022441
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6067633333333333
$ws.Range("H7").Value = 1.82029
$ws.Range("I7").Value = 0.01072372580252768
$ws.Range("J7").Value = 0.01072372580252768
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.7602199999999999
$ws.Range("N7").Value = 2.28066
$ws.Range("O7").Value = 0.5087092361331531
$ws.Range("P7").Value = 0.5087092361331531
$ws.Range("Q7").Value = 0.4612736212666666
$ws.Range("R7").Value = 4.1514625914
$ws.Range("S7").Value = 0.005455258361505242
$ws.Range("T7").Value = 0.005455258361505242
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4787493333333333
$ws.Range("H8").Value = 1.436248
$ws.Range("I8").Value = 0.008461250535040448
$ws.Range("J8").Value = 0.00846125053504045
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7341896666666666
$ws.Range("N8").Value = 2.202569
$ws.Range("O8").Value = 0.4912907638668469
$ws.Range("P8").Value = 0.4912907638668469
$ws.Range("Q8").Value = 0.3514928134568889
$ws.Range("R8").Value = 3.163435321112
$ws.Range("S8").Value = 0.004156934238628789
$ws.Range("T8").Value = 0.00415693423862879
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4787493333333333
$ws.Range("H9").Value = 1.436248
$ws.Range("I9").Value = 0.008461250535040448
$ws.Range("J9").Value = 0.00846125053504045
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7602199999999999
$ws.Range("N9").Value = 2.28066
$ws.Range("O9").Value = 0.5087092361331531
$ws.Range("P9").Value = 0.5087092361331531
$ws.Range("Q9").Value = 0.3639548181866666
$ws.Range("R9").Value = 3.27559336368
$ws.Range("S9").Value = 0.004304316296411659
$ws.Range("T9").Value = 0.00430431629641166
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.840640999999999
$ws.Range("H10").Value = 20.521923
$ws.Range("I10").Value = 0.1208991288160602
$ws.Range("J10").Value = 0.1208991288160602
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.7341896666666666
$ws.Range("N10").Value = 2.202569
$ws.Range("O10").Value = 0.4912907638668469
$ws.Range("P10").Value = 0.4912907638668469
$ws.Range("Q10").Value = 5.022327935576333
$ws.Range("R10").Value = 45.200951420187
$ws.Range("S10").Value = 0.05939662534687856
$ws.Range("T10").Value = 0.05939662534687855
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.840640999999999
$ws.Range("H11").Value = 20.521923
$ws.Range("I11").Value = 0.1208991288160602
$ws.Range("J11").Value = 0.1208991288160602
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.7602199999999999
$ws.Range("N11").Value = 2.28066
$ws.Range("O11").Value = 0.5087092361331531
$ws.Range("P11").Value = 0.5087092361331531
$ws.Range("Q11").Value = 5.200392101019998
$ws.Range("R11").Value = 46.80352890917999
$ws.Range("S11").Value = 0.06150250346918168
$ws.Range("T11").Value = 0.06150250346918168
